$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 233.33333
$ws.Range("I2").Value = 233.33333
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 233.33333
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -120.33333
$ws.Range("N2").ClearContents()

$ws.Range("H5").Value = 22.833334
$ws.Range("I5").Value = 28.857143
$ws.Range("J5").Value = 14.4
$ws.Range("K5").Value = 28.857143
$ws.Range("L5").Value = 14.4
$ws.Range("M5").Value = 86.14285699999999
$ws.Range("N5").Value = -244.4

$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("N18").ClearContents()

$ws.Range("H33").Value = 177.75
$ws.Range("I33").Value = 142.55556
$ws.Range("J33").Value = 283.33334
$ws.Range("K33").Value = 142.55556
$ws.Range("L33").Value = 283.33334
$ws.Range("M33").Value = 86.44443999999999
$ws.Range("N33").Value = -741.33334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19150.688
$ws.Range("I32").Value = 15366.949
$ws.Range("J32").Value = 63798.8
$ws.Range("K32").Value = 15366.949
$ws.Range("L32").Value = 63798.8
$ws.Range("M32").Value = -15079.949
$ws.Range("N32").Value = -64372.8

$ws.Range("H56").Value = 50000
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 50000
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 50000
$ws.Range("N56").Value = -51484

$ws.Range("H74").Value = 12825185
$ws.Range("I74").Value = 17857824
$ws.Range("J74").Value = 14829.818
$ws.Range("K74").Value = 17857824
$ws.Range("L74").Value = 14829.818
$ws.Range("M74").Value = -17856950
$ws.Range("N74").Value = -16577.818

$ws.Range("H77").Value = 12825185
$ws.Range("I77").Value = 17857824
$ws.Range("J77").Value = 14829.818
$ws.Range("K77").Value = 89289120
$ws.Range("L77").Value = 74149.09
$ws.Range("M77").Value = -89284752
$ws.Range("N77").Value = -82885.09

$ws.Range("H97").Value = 835.069
$ws.Range("I97").Value = 782.3333
$ws.Range("J97").Value = 1088.2
$ws.Range("K97").Value = 782.3333
$ws.Range("L97").Value = 1088.2
$ws.Range("M97").Value = -286.3333
$ws.Range("N97").Value = -2080.2

$ws.Range("H122").Value = 1679.24
$ws.Range("I122").Value = 1621.8572
$ws.Range("J122").Value = 1701.5555
$ws.Range("K122").Value = 4865.571599999999
$ws.Range("L122").Value = 5104.666499999999
$ws.Range("M122").Value = -2415.571599999999
$ws.Range("N122").Value = -10004.6665

$ws.Range("H123").Value = 44990
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 44990
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 44990
$ws.Range("N123").Value = -54790

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()

$ws.Range("H132").Value = 1016054.7
$ws.Range("I132").Value = 1328137.2
$ws.Range("J132").Value = 1786.25
$ws.Range("K132").Value = 3984411.6
$ws.Range("L132").Value = 5358.75
$ws.Range("M132").Value = -3981881.6
$ws.Range("N132").Value = -10418.75

$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H128").Value = 13333.333
$ws.Range("I128").Value = 13333.333
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 39999.999
$ws.Range("L128").Value = 0
$ws.Range("M128").Value = -37509.999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5847.4614
$ws.Range("I31").Value = 1937.8334
$ws.Range("J31").Value = 7585.074
$ws.Range("K31").Value = 1937.8334
$ws.Range("L31").Value = 7585.074
$ws.Range("M31").Value = -1642.8334

$ws.Range("H34").Value = 5847.4614
$ws.Range("I34").Value = 1937.8334
$ws.Range("J34").Value = 7585.074
$ws.Range("K34").Value = 1937.8334
$ws.Range("L34").Value = 7585.074
$ws.Range("M34").Value = -1735.8334

$ws.Range("H41").Value = 2566.6667
$ws.Range("I41").Value = 2566.6667
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 2566.6667
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -2138.6667
$ws.Range("N41").ClearContents()

$ws.Range("H50").Value = 10646
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 10646
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 10646
$ws.Range("N50").Value = -11896

$ws.Range("H51").Value = 9898.700000000001
$ws.Range("I51").Value = 8800
$ws.Range("J51").Value = 10173.375
$ws.Range("K51").Value = 8800
$ws.Range("L51").Value = 10173.375
$ws.Range("M51").Value = -8064
$ws.Range("N51").Value = -11645.375

$ws.Range("H59").Value = 15295.3
$ws.Range("I59").Value = 16000
$ws.Range("J59").Value = 15217
$ws.Range("K59").Value = 16000
$ws.Range("L59").Value = 15217
$ws.Range("M59").Value = -14855
$ws.Range("N59").Value = -17507

$ws.Range("H60").Value = 9741.200000000001
$ws.Range("I60").Value = 8000
$ws.Range("J60").Value = 9934.666999999999
$ws.Range("K60").Value = 8000
$ws.Range("L60").Value = 9934.666999999999
$ws.Range("M60").Value = -7489
$ws.Range("N60").Value = -10956.667

$ws.Range("H61").Value = 9898.700000000001
$ws.Range("I61").Value = 8800
$ws.Range("J61").Value = 10173.375
$ws.Range("K61").Value = 8800
$ws.Range("L61").Value = 10173.375
$ws.Range("M61").Value = -8452
$ws.Range("N61").Value = -10869.375

$ws.Range("H68").Value = 18226
$ws.Range("I68").Value = 14268
$ws.Range("J68").Value = 18720.75
$ws.Range("K68").Value = 14268
$ws.Range("L68").Value = 18720.75
$ws.Range("M68").Value = -13519
$ws.Range("N68").Value = -20218.75

$ws.Range("H71").Value = 18226
$ws.Range("I71").Value = 14268
$ws.Range("J71").Value = 18720.75
$ws.Range("K71").Value = 42804
$ws.Range("L71").Value = 56162.25
$ws.Range("M71").Value = -39060
$ws.Range("N71").Value = -63650.25

$ws.Range("H74").Value = 14428.363
$ws.Range("I74").Value = 5185
$ws.Range("J74").Value = 16482.445
$ws.Range("K74").Value = 5185
$ws.Range("L74").Value = 16482.445
$ws.Range("M74").Value = -4311
$ws.Range("N74").Value = -18230.445

$ws.Range("H77").Value = 14428.363
$ws.Range("I77").Value = 5185
$ws.Range("J77").Value = 16482.445
$ws.Range("K77").Value = 15555
$ws.Range("L77").Value = 49447.335
$ws.Range("M77").Value = -11187
$ws.Range("N77").Value = -58183.335

$ws.Range("H99").Value = 126401
$ws.Range("I99").Value = 1270
$ws.Range("J99").Value = 251532
$ws.Range("K99").Value = 1270
$ws.Range("L99").Value = 251532
$ws.Range("M99").Value = 228
$ws.Range("N99").Value = -254528

$ws.Range("H122").Value = 9228.076999999999
$ws.Range("I122").Value = 15994.571
$ws.Range("J122").Value = 1333.8334
$ws.Range("K122").Value = 47983.713
$ws.Range("L122").Value = 4001.5002
$ws.Range("M122").Value = -45533.713
$ws.Range("N122").Value = -8901.5002

$ws.Range("H126").Value = 126401
$ws.Range("I126").Value = 1270
$ws.Range("J126").Value = 251532
$ws.Range("K126").Value = 3810
$ws.Range("L126").Value = 754596
$ws.Range("M126").Value = -1340
$ws.Range("N126").Value = -759536

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 7003950
$ws.Range("I122").Value = 16667294
$ws.Range("J122").Value = 769534.6
$ws.Range("K122").Value = 150005646
$ws.Range("L122").Value = 6925811.399999999
$ws.Range("M122").Value = -150003196
$ws.Range("N122").Value = -6930711.399999999

$ws.Range("H131").Value = 72924024
$ws.Range("I131").Value = 111121224
$ws.Range("J131").Value = 50005704
$ws.Range("K131").Value = 333363672
$ws.Range("L131").Value = 150017112
$ws.Range("M131").Value = -333358632
$ws.Range("N131").Value = -150027192

$ws.Range("H132").Value = 55556756
$ws.Range("I132").Value = 83334310
$ws.Range("J132").Value = 1646.6666
$ws.Range("K132").Value = 750008790
$ws.Range("L132").Value = 14819.9994
$ws.Range("M132").Value = -750006260
$ws.Range("N132").Value = -19879.9994

$ws.Range("H133").Value = 5000
$ws.Range("I133").Value = 5000
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 15000
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -9940
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1311.4286
$ws.Range("I102").Value = 1240
$ws.Range("J102").Value = 1340
$ws.Range("K102").Value = 1240
$ws.Range("L102").Value = 1340
$ws.Range("M102").Value = 382
$ws.Range("N102").Value = -4584

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2067.484
$ws.Range("I68").Value = 1833.7894
$ws.Range("J68").Value = 2437.5
$ws.Range("K68").Value = 1833.7894
$ws.Range("L68").Value = 2437.5
$ws.Range("M68").Value = -1084.7894
$ws.Range("N68").Value = -3935.5

$ws.Range("H71").Value = 2067.484
$ws.Range("I71").Value = 1833.7894
$ws.Range("J71").Value = 2437.5
$ws.Range("K71").Value = 9168.947
$ws.Range("L71").Value = 12187.5
$ws.Range("M71").Value = -5424.947
$ws.Range("N71").Value = -19675.5

$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
